$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "6am - 2pm"
$ws.Range("A2").Value = "2pm - 10am"

$ws.Range("C9").Select()
